$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.548.32"
$ws.Range("E2").Value = "  +0.40%  "

$ws.Range("D3").Value = "1.816.66"
$ws.Range("E3").Value = "  +0.69%  "

$ws.Range("E4").Value = "  -0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "225.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.03%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.606"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.86%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +0.03%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "46.02"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +26.72%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.297"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +1.55%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0677"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.19%  "

$ws.Range("E11").Value = "  +3.48%  "

$ws.Range("D12").Value = "2.081.64"
$ws.Range("E12").Value = "  +0.90%  "

$ws.Range("D13").Value = "1.823.08"
$ws.Range("E13").Value = "  +1.07%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "11.15"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.18%  "

$ws.Range("E15").Value = "  +1.99%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "4.49"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.28%  "

$ws.Range("D17").Value = "34.575.41"
$ws.Range("E17").Value = "  +0.53%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.07"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.76%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "241.99"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -0.26%  "

$ws.Range("E20").Value = "  +1.06%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.70"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.17%  "

$ws.Range("E22").Value = "  -0.01%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "4.49"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +9.71%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "2.18"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.03%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "170.75"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.30%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.82"
$ws.Range("D26").Style = "Normal"

$ws.Range("E27").Value = "  +1.98%  "

$ws.Range("E28").Value = "  +0.27%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.997"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.32%  "

$ws.Range("E30").Value = "  +1.73%  "

$ws.Range("E31").Value = "  +1.21%  "

$ws.Range("E32").Value = "  +0.39%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0521"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +1.23%  "

$ws.Range("E34").Value = "  +2.95%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "89.77"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +11.30%  "

$ws.Range("E36").Value = "  +1.28%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "15.48"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +16.40%  "

$ws.Range("B38").Value = "Maker"
$ws.Range("C38").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D38").Value = "1.323.58"
$ws.Range("E38").Value = "  -2.80%  "

$ws.Range("B39").Value = "TrustWalletToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "1.06"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.30%  "

$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.40"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.73%  "

$ws.Range("E41").Value = "  +2.78%  "

$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.960"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +2.36%  "

$ws.Range("B43").Value = "WEMIXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.22"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +5.15%  "

$ws.Range("E44").Value = "  +1.27%  "

$ws.Range("E45").Value = "  -0.14%  "

$ws.Range("E46").Value = "  +3.85%  "

$ws.Range("D47").Value = "1.984.01"
$ws.Range("E47").Value = "  +1.04%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "5.91"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.29%  "

$ws.Range("E49").Value = "  -0.15%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "101.09"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.32%  "

$ws.Range("E51").Value = "  +0.79%  "
